# The commit swaps the data between row 3 ("Minnesota Supreme Court to hear
# Byron Smith's appeal") and row 4 ("Crime scene photos released from Byron
# Smith murder trial (PHOTOS)") so that the "Crime scene photos" story now
# appears before the "Minnesota Supreme Court" story (sorted by timestamp -
# this is effectively inserting a new time-bucket-analysis JSON entry that
# reorders these two rows by their timestamp).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture current row 3 values (title, timestamp, historical distance, time bucket, uri)
$row3Title = $ws.Range("A3").Value2
$row3Timestamp = $ws.Range("B3").Value2
$row3Distance = $ws.Range("C3").Value2
$row3Bucket = $ws.Range("D3").Value2
$row3Uri = $ws.Range("E3").Value2

# Capture current row 4 values
$row4Title = $ws.Range("A4").Value2
$row4Timestamp = $ws.Range("B4").Value2
$row4Distance = $ws.Range("C4").Value2
$row4Bucket = $ws.Range("D4").Value2
$row4Uri = $ws.Range("E4").Value2

# Write row 4's content into row 3
$ws.Range("A3").Value = $row4Title
$ws.Range("B3").Value = $row4Timestamp
$ws.Range("C3").Value = $row4Distance
$ws.Range("D3").Value = $row4Bucket
$ws.Range("E3").Value = $row4Uri

# Write row 3's original content into row 4
$ws.Range("A4").Value = $row3Title
$ws.Range("B4").Value = $row3Timestamp
$ws.Range("C4").Value = $row3Distance
$ws.Range("D4").Value = $row3Bucket
$ws.Range("E4").Value = $row3Uri

# The hyperlinks on E2:E5 must follow the swapped uri text. Remove every
# hyperlink on the sheet and re-add them in row order so each E-cell's
# hyperlink target matches its (now swapped) displayed uri text.
$ws.Cells.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("E2"), $ws.Range("E2").Value2) | Out-Null
$ws.Hyperlinks.Add($ws.Range("E3"), $ws.Range("E3").Value2) | Out-Null
$ws.Hyperlinks.Add($ws.Range("E4"), $ws.Range("E4").Value2) | Out-Null
$ws.Hyperlinks.Add($ws.Range("E5"), $ws.Range("E5").Value2) | Out-Null

# Hyperlinks.Add() re-derives a fresh cell style; restore the sheet's
# original named "Hyperlink" style on the uri column so the cell formatting
# matches what it was before (only the underlying link target changed).
$ws.Range("E2:E5").Style = "Hyperlink"
